$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting for numeric-looking price cells so Excel
# does not reinterpret them as numbers (they are stored as text).
$textCells = @("D5", "D6", "D7", "D8", "D10", "D11", "D16", "D17", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D28", "D30", "D33", "D34", "D36", "D37", "D39", "D40", "D41", "D42", "D43", "D48", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '54.375.17'
$ws.Range("E2").Value = '  -6.97%  '
$ws.Range("D3").Value = '2.427.80'
$ws.Range("E3").Value = '  -10.75%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '464.10'
$ws.Range("E5").Value = '  -7.64%  '
$ws.Range("D6").Value = '130.92'
$ws.Range("E6").Value = '  -6.87%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.27%  '
$ws.Range("D8").Value = '0.491'
$ws.Range("E8").Value = '  -7.31%  '
$ws.Range("D9").Value = '2.426.17'
$ws.Range("E9").Value = '  -11.34%  '
$ws.Range("D10").Value = '0.0947'
$ws.Range("E10").Value = '  -9.61%  '
$ws.Range("D11").Value = '5.30'
$ws.Range("E11").Value = '  -12.38%  '
$ws.Range("E12").Value = '  -9.77%  '
$ws.Range("E13").Value = '  -4.24%  '
$ws.Range("D14").Value = '2.853.78'
$ws.Range("E14").Value = '  -10.82%  '
$ws.Range("D15").Value = '54.400.82'
$ws.Range("E15").Value = '  -7.14%  '
$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D16").Value = '19.63'
$ws.Range("E16").Value = '  -9.48%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = '0.0000131'
$ws.Range("E17").Value = '  -2.72%  '
$ws.Range("D18").Value = '2.431.68'
$ws.Range("E18").Value = '  -11.05%  '
$ws.Range("D19").Value = '4.19'
$ws.Range("E19").Value = '  -12.02%  '
$ws.Range("D20").Value = '309.43'
$ws.Range("E20").Value = '  -9.73%  '
$ws.Range("D21").Value = '9.50'
$ws.Range("E21").Value = '  -13.44%  '
$ws.Range("D22").Value = '0.996'
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("D23").Value = '5.67'
$ws.Range("E23").Value = '  +0.22%  '
$ws.Range("D24").Value = '5.35'
$ws.Range("E24").Value = '  -14.48%  '
$ws.Range("D25").Value = '56.20'
$ws.Range("E25").Value = '  -10.69%  '
$ws.Range("E26").Value = '  +1.42%  '
$ws.Range("E27").Value = '  -10.06%  '
$ws.Range("D28").Value = '0.155'
$ws.Range("E28").Value = '  -10.02%  '
$ws.Range("D29").Value = '2.533.35'
$ws.Range("E29").Value = '  -10.85%  '
$ws.Range("D30").Value = '7.10'
$ws.Range("E30").Value = '  -5.28%  '
$ws.Range("E31").Value = '  +0.23%  '
$ws.Range("D32").Value = '0.0₃0708'
$ws.Range("E32").Value = '  -14.46%  '
$ws.Range("D33").Value = '146.42'
$ws.Range("E33").Value = '  -3.22%  '
$ws.Range("D34").Value = '17.66'
$ws.Range("E34").Value = '  -7.87%  '
$ws.Range("E35").Value = '  -10.69%  '
$ws.Range("D36").Value = '4.98'
$ws.Range("E36").Value = '  -8.35%  '
$ws.Range("D37").Value = '3.52'
$ws.Range("E37").Value = '  -15.94%  '
$ws.Range("E38").Value = '  -7.40%  '
$ws.Range("D39").Value = '0.796'
$ws.Range("E39").Value = '  -16.09%  '
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  +0.16%  '
$ws.Range("D41").Value = '32.90'
$ws.Range("E41").Value = '  -8.41%  '
$ws.Range("D42").Value = '0.593'
$ws.Range("E42").Value = '  -1.39%  '
$ws.Range("D43").Value = '0.0520'
$ws.Range("E43").Value = '  -6.80%  '
$ws.Range("E44").Value = '  -8.94%  '
$ws.Range("E45").Value = '  -2.94%  '
$ws.Range("E46").Value = '  -11.87%  '
$ws.Range("D47").Value = '1.926.87'
$ws.Range("E47").Value = '  -11.89%  '
$ws.Range("D48").Value = '0.0875'
$ws.Range("E48").Value = '  -1.14%  '
$ws.Range("E49").Value = '  -4.66%  '
$ws.Range("D50").Value = '16.41'
$ws.Range("E50").Value = '  -13.76%  '
$ws.Range("B51").Value = 'Bittensor'
$ws.Range("C51").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D51").Value = '227.83'
$ws.Range("E51").Value = '  +0.56%  '
